$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-14 03:15:02"
$wsZhCn.Range("G2").Value = "2016-01-14 03:15:50"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-14 03:15:14"
$wsDeDe.Range("G2").Value = "2016-01-14 03:16:11"
